$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BD2 151 -> 126
$ws.Range("BD2").Value = 126

# Row 3: G3, I3, N3, AQ3, AZ3
$ws.Range("G3").Value = 2.9
$ws.Range("I3").Value = 2.88
$ws.Range("N3").Value = 4.75
$ws.Range("AQ3").Value = 81
$ws.Range("AZ3").Value = 67

# Row 5: M5, N5
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5

# Row 8: Q8, R8
$ws.Range("Q8").Value = 2.4
$ws.Range("R8").Value = 1.53

# Delete row 10 entirely (the old Uruguay row, previously row 11, shifts up to row 10)
$ws.Rows(10).Delete()
